# Apply updated TPM values to the LR-pairs sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Ccl5/ECs -> Ccr1/MuSCs) updated values
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.021814
$ws.Range("H2").Value = 0.065442
$ws.Range("I2").Value = 0.1008129179549036
$ws.Range("J2").Value = 0.1008129179549036
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.009232000000000001
$ws.Range("N2").Value = 0.027696
$ws.Range("Q2").Value = 0.000201386848
$ws.Range("R2").Value = 0.001812481632
$ws.Range("S2").Value = 0.1008129179549036
$ws.Range("T2").Value = 0.1008129179549036

# Row 3 (Ccl5/MuSCs -> Ccr1/MuSCs) updated values
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.194567
$ws.Range("H3").Value = 0.583701
$ws.Range("I3").Value = 0.8991870820450963
$ws.Range("J3").Value = 0.8991870820450963
$ws.Range("K3").Value = 2
$ws.Range("L3").Value = 0.6666666666666666
$ws.Range("M3").Value = 0.009232000000000001
$ws.Range("N3").Value = 0.027696
$ws.Range("Q3").Value = 0.001796242544
$ws.Range("R3").Value = 0.016166182896
$ws.Range("S3").Value = 0.8991870820450963
$ws.Range("T3").Value = 0.8991870820450963
